$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update aggregate metrics to reflect the new closed trade
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1500.77   # Current Capital
$summary.Range("B4").Value = 0.77      # Total P&L $
$summary.Range("B5").Value = 0.64      # Total P&L %
$summary.Range("B6").Value = 24        # Total Trades
$summary.Range("B8").Value = 9         # Losing Trades
$summary.Range("B9").Value = 54.17     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking strategy row (row 6)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 100.77     # Capital
$status.Range("D6").Value = 24         # Trades
$status.Range("E6").Value = 0.77       # P&L $
$status.Range("F6").Value = 0.77       # P&L %
$status.Range("G6").Value = 54.17      # Win Rate %

# ---------------------------------------------------------------------------
# Append the newly closed trade (#24) as row 25 to both the "All Trades"
# sheet and the per-strategy "MarketMaking" sheet.
# ---------------------------------------------------------------------------
$sheetsToAppend = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetsToAppend) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Force column B (Date) to be stored as plain text, matching the
    # existing rows, instead of being auto-converted to a date serial.
    $ws.Range("B25").NumberFormat = "@"

    $ws.Range("A25").Value = 24
    $ws.Range("B25").Value = "2026-02-17"
    $ws.Range("C25").Value = "23:57:25"
    $ws.Range("D25").Value = "MarketMaking"
    $ws.Range("E25").Value = "DOWN"
    $ws.Range("F25").Value = 0.36
    $ws.Range("G25").Value = 0.3
    $ws.Range("H25").Value = "CLOSED"
    $ws.Range("I25").Value = -16.6667
    $ws.Range("J25").Value = -0.06
    $ws.Range("K25").Value = 100.77
    $ws.Range("L25").Value = 0
    $ws.Range("M25").Value = 0
    $ws.Range("N25").Value = 0.6
    $ws.Range("O25").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P25").Value = "early_exit"
    $ws.Range("Q25").Value = 0.14
}
